$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

Set-TextValue $ws 'D2' '320.00'
Set-TextValue $ws 'E2' '4.71%'
Set-TextValue $ws 'D3' '36.10'
Set-TextValue $ws 'E3' '-0.07%'
Set-TextValue $ws 'E4' '1.05%'
Set-TextValue $ws 'D5' '0.08173'
Set-TextValue $ws 'E5' '4.10%'
Set-TextValue $ws 'D6' '2.148'
Set-TextValue $ws 'E6' '-2.31%'
Set-TextValue $ws 'D7' '8.039'
Set-TextValue $ws 'E7' '1.44%'
$ws.Range('B8').Value = 'MXToken'
$ws.Range('C8').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
Set-TextValue $ws 'D8' '0.9267'
Set-TextValue $ws 'E8' '0.67%'
$ws.Range('B9').Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range('C9').Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
Set-TextValue $ws 'D9' '0.1008'
Set-TextValue $ws 'E9' '4.82%'
$ws.Range('B10').Value = 'WazirX'
$ws.Range('C10').Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
Set-TextValue $ws 'D10' '0.1884'
Set-TextValue $ws 'E10' '1.22%'
$ws.Range('B11').Value = 'MandalaExchangeToken'
$ws.Range('C11').Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
Set-TextValue $ws 'D11' '0.09220'
Set-TextValue $ws 'E11' '6.76%'
$ws.Range('B12').Value = 'BitrueCoin'
$ws.Range('C12').Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
Set-TextValue $ws 'D12' '0.03595'
Set-TextValue $ws 'E12' '3.38%'
$ws.Range('B13').Value = 'BitMartToken'
$ws.Range('C13').Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
Set-TextValue $ws 'D13' '0.09925'
Set-TextValue $ws 'E13' '-0.11%'
$ws.Range('B14').Value = 'BitForexToken'
$ws.Range('C14').Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
Set-TextValue $ws 'D14' '0.001432'
Set-TextValue $ws 'E14' '0.17%'
$ws.Range('B15').Value = 'TigerCash'
$ws.Range('C15').Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
Set-TextValue $ws 'D15' '0.005694'
Set-TextValue $ws 'E15' '1.16%'
$ws.Range('B16').Value = 'LEO'
$ws.Range('C16').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
Set-TextValue $ws 'D16' '3.466'
Set-TextValue $ws 'E16' '0.02%'
$ws.Range('B17').Value = 'GateToken'
$ws.Range('C17').Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
Set-TextValue $ws 'D17' '4.135'
Set-TextValue $ws 'E17' '0.95%'
Set-TextValue $ws 'E18' '16.93%'
Set-TextValue $ws 'E19' '-1.57%'
Set-TextValue $ws 'D20' '0.1328'
Set-TextValue $ws 'E20' '2.11%'
Set-TextValue $ws 'D21' '5.071'
Set-TextValue $ws 'E21' '5.23%'
Set-TextValue $ws 'D22' '0.2189'
Set-TextValue $ws 'E22' '-0.81%'
Set-TextValue $ws 'D23' '0.04599'
Set-TextValue $ws 'E23' '1.40%'
Set-TextValue $ws 'D24' '0.001243'
Set-TextValue $ws 'E24' '0.71%'
Set-TextValue $ws 'D25' '0.004727'
Set-TextValue $ws 'E25' '-7.22%'
Set-TextValue $ws 'D26' '0.0001301'
Set-TextValue $ws 'E26' '-7.16%'
Set-TextValue $ws 'D27' '0.0004502'
Set-TextValue $ws 'E27' '-5.28%'
Set-TextValue $ws 'D39' '0.02015'
Set-TextValue $ws 'E39' '10.28%'
Set-TextValue $ws 'D40' '0.04996'
Set-TextValue $ws 'E40' '4.52%'
Set-TextValue $ws 'D41' '0.007783'
Set-TextValue $ws 'E41' '0.92%'
Set-TextValue $ws 'D42' '0.1402'
Set-TextValue $ws 'E42' '0.32%'
Set-TextValue $ws 'D43' '0.007817'
Set-TextValue $ws 'E43' '0.97%'
Set-TextValue $ws 'D44' '0.002097'
Set-TextValue $ws 'E44' '-6.05%'
Set-TextValue $ws 'D45' '0.01212'
Set-TextValue $ws 'E45' '8.14%'
Set-TextValue $ws 'D46' '0.00006485'
Set-TextValue $ws 'E46' '3.60%'
Set-TextValue $ws 'E47' '-0.03%'
Set-TextValue $ws 'D49' '0.001901'
Set-TextValue $ws 'E49' '-5.01%'
Set-TextValue $ws 'D50' '0.00002101'
Set-TextValue $ws 'E50' '-0.03%'
Set-TextValue $ws 'D51' '0.0002001'
Set-TextValue $ws 'E51' '-0.03%'
